# Applies crypto price/volume updates for Mon Aug 14 08:32:41 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.392.24"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.847.14"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.34"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6267"
$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07500"
$ws.Range("E8").Value = "  -1.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2904"
$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.45"
$ws.Range("E10").Value = "  -1.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07736"
$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("D12").Value = "1.846.48"
$ws.Range("E12").Value = "  -2.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.999"
$ws.Range("E13").Value = "  -0.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6805"
$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001059"
$ws.Range("E15").Value = "  +0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.21"
$ws.Range("E16").Value = "  -1.29%  "

$ws.Range("D17").Value = "2.106.52"
$ws.Range("E17").Value = "  -3.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.169"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("D19").Value = "29.430.22"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.51"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.33"
$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.482"
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9995"
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.17"
$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.410"
$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.53"
$ws.Range("E28").Value = "  -0.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06442"
$ws.Range("E29").Value = "  +15.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.413"
$ws.Range("E30").Value = "  +2.12%  "

$ws.Range("E31").Value = "  +1.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.093"
$ws.Range("E32").Value = "  -0.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.097"
$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("E35").Value = "  -1.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6966"
$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.578"
$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("D38").Value = "1.270.70"
$ws.Range("E38").Value = "  +3.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01836"
$ws.Range("E39").Value = "  +1.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.824"
$ws.Range("E40").Value = "  +4.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.720"
$ws.Range("E41").Value = "  +5.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9090"
$ws.Range("E42").Value = "  +0.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9989"
$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("D44").Value = "2.010.24"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.31"
$ws.Range("E45").Value = "  -0.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.31"
$ws.Range("E46").Value = "  +0.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.742"
$ws.Range("E47").Value = "  +3.67%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.078"
$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("E49").Value = "  +3.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.102"
$ws.Range("E50").Value = "  +1.31%  "

# Row 51: coin swapped from BabyDogeCoin to TheSandbox
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3957"
$ws.Range("E51").Value = "  -1.39%  "
